$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Collapse the split runs in the title / author / abstract paragraphs into
#    single runs (the wording itself is unchanged, only the run segmentation
#    is merged).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Questions: Factorization", $false, $false, $false, $false, $false,
    $true, 1, $false, "Questions: Factorization", 2) | Out-Null

$d.Content.Find.Execute(
    "Millie Pike", $false, $false, $false, $false, $false,
    $true, 1, $false, "Millie Pike", 2) | Out-Null

$d.Content.Find.Execute(
    "A selection of questions for the study guide on factorization.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "A selection of questions for the study guide on factorization.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. In the two equations that render as  x(x-6)+3(6-x)[=0], swap the order
#    of <m:sepChr/> and <m:endChr/> inside each delimiter's <m:dPr/> so it
#    becomes begChr, sepChr, endChr, grow (was begChr, endChr, sepChr, grow).
#    This has no visual/semantic effect - it only reorders markup - so we
#    rebuild each affected <m:oMath/> with the swapped order and splice it
#    back in via InsertXML (WordOpenXML is read-only on Range).
# ---------------------------------------------------------------------------
$oldOrder = '<m:begChr m:val="(" /><m:endChr m:val=")" /><m:sepChr m:val="" />'
$newOrder = '<m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" />'

$oMaths = $d.OMaths
for ($i = 1; $i -le $oMaths.Count; $i++) {
    $om = $oMaths.Item($i)
    $currentXml = $om.Range.WordOpenXML
    if ($currentXml -match "sepChr") {
        $match = [regex]::Match($currentXml, "<m:oMath[\s\S]*?</m:oMath>")
        $fixedFragment = $match.Value.Replace($oldOrder, $newOrder)

        $package = '<?xml version="1.0" standalone="yes"?>' +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math">' +
            '<w:body><w:p>' + $fixedFragment + '</w:p></w:body></w:document>' +
            '</pkg:xmlData></pkg:part></pkg:package>'

        $om.Range.InsertXML($package)
    }
}

# The WordOpenXML round trip above renders the math minus sign (U+2212) as a
# plain hyphen-minus; restore the original glyph in the two equations that
# were touched so only the dPr element order actually changed.
$d.Content.Find.Execute(
    "x(x-6)+3(6-x)=0", $false, $false, $false, $false, $false,
    $true, 1, $false, "x(x`u{2212}6)+3(6`u{2212}x)=0", 2) | Out-Null

$d.Content.Find.Execute(
    "x(x-6)+3(6-x)", $false, $false, $false, $false, $false,
    $true, 1, $false, "x(x`u{2212}6)+3(6`u{2212}x)", 2) | Out-Null
